# Add a new "Rounding_thresholds" sheet at the end of the workbook.
# It lists every model parameter together with a small-value threshold
# (below which a value is treated as "effectively zero") and the value
# that should replace it.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Rounding_thresholds"

# Header row
$ws.Range("A1").Value = "Parameter"
$ws.Range("B1").Value = "Threshold"
$ws.Range("C1").Value = "Replace with"
$ws.Range("A1:C1").Font.Bold = $true

$parameters = @(
    "Par_AnnualEmissionLimit",
    "Par_AnnualExogenousEmission",
    "Par_AnnualMaxNewCapacity",
    "Par_AnnualMinNewCapacity",
    "Par_AnnualSectoralEmissionLimit",
    "Par_AvailabilityFactor",
    "Par_BaseYearProduction",
    "Par_BaseYearSlack",
    "Par_CapacityFactor",
    "Par_CapacityToActivityUnit",
    "Par_CapitalCost",
    "Par_CapitalCostStorage",
    "Par_CommissionedTradeCapacity",
    "Par_EmissionActivityRatio",
    "Par_EmissionContentPerFuel",
    "Par_EmissionPenaltyTagTech",
    "Par_EmissionsPenalty",
    "Par_FixedCost",
    "Par_GeneralDiscountRate",
    "Par_GrowthRateTradeCapacity",
    "Par_InputActivityRatio",
    "Par_MinStorageCharge",
    "Par_ModalSplitByFuel",
    "Par_ModelPeriodActivityMaxLimit",
    "Par_ModelPeriodEmissionLimit",
    "Par_ModelPeriodExogenousEmission",
    "Par_NewCapacityExpansionStop",
    "Par_OperationalLife",
    "Par_OperationalLifeStorage",
    "Par_OutputActivityRatio",
    "Par_ProductionChangeCost",
    "Par_ProductionGrowthLimit",
    "Par_REMinProductionTarget",
    "Par_RETagFuel",
    "Par_RETagTechnology",
    "Par_RampingDownFactor",
    "Par_RampingUpFactor",
    "Par_RegionalAnnualEmissionLimit",
    "Par_RegionalBaseYearProduction",
    "Par_RegionalCCSLimit",
    "Par_RegionalModelPeriodEmission",
    "Par_ReserveMargin",
    "Par_ReserveMarginTagFuel",
    "Par_ReserveMarginTagTechnology",
    "Par_ResidualCapacity",
    "Par_ResidualStorageCapacity",
    "Par_SelfSufficiency",
    "Par_SocialDiscountRate",
    "Par_SpecifiedAnnualDemand",
    "Par_SpecifiedDemandDevelopment",
    "Par_StorageE2PRatio",
    "Par_StorageLevelStart",
    "Par_TagCanFuelBeTraded",
    "Par_TagDemandFuelToSector",
    "Par_TagElectricTechnology",
    "Par_TagFuelToSubsets",
    "Par_TagModalTypeToModalGroups",
    "Par_TagTechnologyToModalType",
    "Par_TagTechnologyToSector",
    "Par_TagTechnologyToSubsets",
    "Par_TagTimeIndependentFuel",
    "Par_TechnologyDiscountRate",
    "Par_TechnologyFromStorage",
    "Par_TechnologyToStorage",
    "Par_TotalAnnualMaxActivity",
    "Par_TotalAnnualMaxCapacity",
    "Par_TotalAnnualMinActivity",
    "Par_TotalAnnualMinCapacity",
    "Par_TradeCapacity",
    "Par_TradeCapacityGrowthCosts",
    "Par_TradeCostsFactor",
    "Par_TradeLossFactor",
    "Par_TradeRoute",
    "Par_VariableCost"
)

$row = 2
foreach ($p in $parameters) {
    $ws.Cells.Item($row, 1).Value = $p
    $ws.Cells.Item($row, 2).Value = 0.00001
    $ws.Cells.Item($row, 2).NumberFormat = "0.00E+00"
    $ws.Cells.Item($row, 3).Value = 0
    $row++
}

$ws.Range("D72").Select()
$ws.Activate()
